# Update the "Ignores or opposes needs" diamond-diagram caption on slide 1
# to read "Lacks community" / "involvement" (two centred lines), per the
# commit "Update to Lacks Community Involvement".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The caption lives inside a top-level group ("Group 74") as the shape
# "Rectangle 76". Search for it by walking the shapes/groups instead of
# relying on fixed indices, so the script is resilient to shape reordering.

function Find-ShapeByText {
    # Positional parameters -- this interpreter does not reliably bind
    # PowerShell's named ("-paramName value") call syntax for functions.
    param($shapes, $needle)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)

        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -like "*$needle*") {
                return $shp
            }
        }

        if ($shp.Type -eq 6) {
            # msoGroup
            $found = Find-ShapeByText $shp.GroupItems $needle
            if ($found) {
                return $found
            }
        }
    }

    return $null
}

$target = Find-ShapeByText $s.Shapes "Ignores or opposes needs"

if ($target) {
    # Two centred paragraphs, split with a carriage return so PowerPoint
    # creates a new <a:p> (preserving the existing run/paragraph formatting)
    # rather than a soft line break within a single paragraph.
    $target.TextFrame.TextRange.Text = "Lacks community`rinvolvement"
} else {
    throw "Could not locate the 'Ignores or opposes needs' caption shape on slide 1"
}
